# Fixing Bulk Operation Template
#
# The "Create Item Category" bulk-upload template (2 columns: Item
# Category Name / Item Group Name) is repurposed into the "Update Item"
# bulk-upload template (6 columns: Item Code, Item Name, Category, Group,
# Unit, Sell Price).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet identity -------------------------------------------------
$ws.Name = "Update Item"

# --- Header row -------------------------------------------------------
# Wipe the old 2-column header and lay out the new 6-column header.
$ws.Range("A1:B1").ClearContents()

$ws.Range("A1").Value = "Item Code"
$ws.Range("B1").Value = "Item Name"
$ws.Range("C1").Value = "Category"
$ws.Range("D1").Value = "Group"
$ws.Range("E1").Value = "Unit"
$ws.Range("F1").Value = "Sell Price"

# A1 keeps its original bold / filled "header" look. B1:F1 get the
# lighter fill-only header look (same fill color as A1, default font) -
# first strip any inherited formatting from the old B1 so the new fill
# doesn't pick up the bold header font along with it.
$ws.Range("B1:F1").ClearFormats()
$ws.Range("B1:F1").Interior.Color = $ws.Range("A1").Interior.Color
$ws.Range("B1:F1").Interior.Pattern = $ws.Range("A1").Interior.Pattern

# --- Column sizing ------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.42
$ws.Columns.Item(2).ColumnWidth = 40.42
$ws.Columns.Item(3).ColumnWidth = 14.92
$ws.Columns.Item(4).ColumnWidth = 12.09
$ws.Columns.Item(5).ColumnWidth = 10.42
$ws.Columns.Item(6).ColumnWidth = 14.75

# --- Leave a clean A1 selection instead of the stale D3 one from the
# old template.
[void]$ws.Range("A1").Select()
